# Applies the "Actualización automática del dashboard" update:
#  1. Fixes H column (fecha_comentario) cells that were left with the wrong
#     date-only number format (numFmtId 167) so they match the other date
#     cells in the column (numFmtId 165, "YYYY-MM-DD HH:MM:SS").
#  2. Fixes J column (likes_count) cells that were written as text ("0"/"6")
#     instead of numbers, converting them to real numeric values.
#  3. Restores the N column (created_time_raw) JSON payloads for rows 9/10,
#     which had been swapped between the two duplicate-post comment rows.
#  4. Adds a new "Resumen_Posts" summary sheet with post-level aggregates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentarios")

# --- 1 & 2: fix H (date format) and J (numeric likes_count) on the rows
#     that were saved with inconsistent types/styles ---
$rowsToFix = @(3, 5, 7, 8, 10, 12, 13, 14, 16, 18)
$likesValues = @{ 3 = 0; 5 = 0; 7 = 0; 8 = 0; 10 = 6; 12 = 0; 13 = 0; 14 = 0; 16 = 0; 18 = 0 }

foreach ($r in $rowsToFix) {
    # Re-apply the same date/time number format used by the correctly
    # formatted rows so Excel reuses that existing style instead of the
    # stray "date only" one.
    $ws.Range("H$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Re-write likes_count as an actual number rather than a text string.
    $ws.Range("J$r").Value = $likesValues[$r]
}

# --- 3: rows 9 and 10 share the same Facebook post/comment but had their
#     created_time_raw (N column) JSON blobs swapped between them; put
#     each one back with the other's value. ---
$n9 = $ws.Range("N9").Value2
$n10 = $ws.Range("N10").Value2
$ws.Range("N9").Value = $n10
$ws.Range("N10").Value = $n9

# --- 4: add the Resumen_Posts summary sheet right after "Comentarios" ---
$summary = $wb.Worksheets.Add($null, $ws)
$summary.Name = "Resumen_Posts"

$summary.Range("A1").Value = "post_number"
$summary.Range("B1").Value = "platform"
$summary.Range("C1").Value = "post_url"
$summary.Range("D1").Value = "Total_Comentarios"
$summary.Range("E1").Value = "Total_Likes"

# Match the bold/centered header style already used on row 1 of Comentarios.
$ws.Range("A1").Copy()
$summary.Range("A1:E1").PasteSpecial(-4122)

$summary.Range("A2").Value = 2
$summary.Range("B2").Value = "Facebook"
$summary.Range("C2").Value = "https://www.facebook.com/100064867445065/posts/1260798862759017/?dco_ad_token=AaprzcNowYg9Z8x7VPPLTUDn0JLBBhvVLiYoWBNA4nuLlQiDP5dy-AXgevLL_V3gTKAoS5-zMd5T54oY&dco_ad_id=120234998089620781"
$summary.Range("D2").Value = 17
$summary.Range("E2").Value = 12

# Keep "Comentarios" as the active sheet/tab, as in the original workbook.
$ws.Activate()

Write-Output "done"
